# Update countries & provincias Spain
# Applies the 27-Sep-2020 03:24 -> 04:41 COVID data refresh to the "Pais" sheet:
#  - Nepal / Venezuela swap rank (row 54/55) as Venezuela's totals overtook Nepal's
#  - Santa Lucia / Timor Oriental swap rank (row 206/207)
#  - Updated case figures for Bolivia, Kazajistan, Nepal, Venezuela, Australia,
#    Islas Turcas y Caicos and San Martin (Parte Holandesa)
#  - "last updated" timestamp bumped from 03:24 to 04:41

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Septiembre de 2020 a las 04:41"

# --- Bolivia (row 31) ---------------------------------------------------
$ws.Range("B31").Value = 133592
$ws.Range("C31").Value = 370
$ws.Range("D31").Value = 93406
$ws.Range("E31").Value = 32358
$ws.Range("G31").Value = 28
$ws.Range("H31").Value = 7828

# --- Kazajistan (row 39) -------------------------------------------------
$ws.Range("B39").Value = 107723
$ws.Range("C39").Value = 64
$ws.Range("E39").Value = 3494

# --- Nepal / Venezuela swap ranking (rows 54-55) --------------------
# Venezuela now has more total cases than Nepal, so it moves up to row 54
# and takes on fresh figures; Nepal drops to row 55 keeping its old figures.
$ws.Range("A54").Value = "Venezuela"
$ws.Range("B54").Value = 71940
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 61528
$ws.Range("E54").Value = 9812
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 600

$ws.Range("A55").Value = "Nepal"
$ws.Range("B55").Value = 71821
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 53013
$ws.Range("E55").Value = 18341
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 467

# --- Australia (row 78) ---------------------------------------------
$ws.Range("B78").Value = 27033
$ws.Range("C78").Value = 17
$ws.Range("E78").Value = 1590
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 872

# --- Islas Turcas y Caicos (row 172) ---------------------------------
$ws.Range("B172").Value = 681
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 620
$ws.Range("E172").Value = 56

# --- San Martin (Parte Holandesa) (row 173) ---------------------------
$ws.Range("B173").Value = 633
$ws.Range("C173").Value = 6
$ws.Range("E173").Value = 79

# --- Santa Lucia / Timor Oriental swap ranking (rows 206-207) --------
$ws.Range("A206").Value = "Santa Lucia"
$ws.Range("A207").Value = "Timor Oriental"
